$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 new values
$ws.Range("A15").Value = 112079417
$ws.Range("B15").Value = 8377
$ws.Range("D15").Value = "LC"
$ws.Range("E15").Value = 106545
$ws.Range("F15").Value = "Mindre märgborre"
$ws.Range("G15").Value = "Tomicus minor"
$ws.Range("H15").Value = "(Hartig, 1834)"
$ws.Range("Q15").Value = 563452
$ws.Range("R15").Value = 6576051

# Row 16 new values
$ws.Range("A16").Value = 112079439
$ws.Range("B16").Value = 90837
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 5966
$ws.Range("F16").Value = "Motaggsvamp"
$ws.Range("G16").Value = "Sarcodon squamosus"
$ws.Range("H16").Value = "(Schaeff.) Quél."
$ws.Range("Q16").Value = 563408
$ws.Range("R16").Value = 6576469
